$d = $word.ActiveDocument

# 1. Replace the text of the second paragraph ("And yet again later on 5/27")
#    with the new text ("And an edit by stanneumann2").
$p2 = $d.Paragraphs(2)
$p2.Range.Find.Execute("And yet again later on 5/27", $false, $false, $false, $false, $false, $true, 1, $false, "And an edit by stanneumann2", 2)

# 2. Apply red (FF0000), 16pt (sz 32 half-points) character formatting to the
#    whole paragraph, including its paragraph mark, so the formatting lands
#    both on the run and on the paragraph's rPr (inside pPr).
$p2 = $d.Paragraphs(2)
$p2.Range.Font.Color = 255
$p2.Range.Font.Size = 16

# 3. Move the "_GoBack" bookmark from the second paragraph to the end of the
#    (empty) paragraph that immediately follows the "Protocol" heading.
$d.Bookmarks("_GoBack").Delete()

$targetPara = $d.Paragraphs(4)
$d.Bookmarks.Add("_GoBack", $targetPara.Range)
